$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 439678.88
$ws.Range("I28").Value = 714831
$ws.Range("J28").Value = 11664.444
$ws.Range("K28").Value = 714831
$ws.Range("L28").Value = 11664.444
$ws.Range("M28").Value = -714346
$ws.Range("N28").Value = -12634.444
$ws.Range("H33").Value = 363.10715
$ws.Range("I33").Value = 285.9091
$ws.Range("J33").Value = 646.1667
$ws.Range("K33").Value = 285.9091
$ws.Range("L33").Value = 646.1667
$ws.Range("M33").Value = -56.90910000000002
$ws.Range("N33").Value = -1104.1667
$ws.Range("H39").Value = 2297.3333
$ws.Range("J39").Value = 2935.75
$ws.Range("L39").Value = 8807.25
$ws.Range("N39").Value = -9399.25
$ws.Range("H41").Value = 305
$ws.Range("I41").Value = 287.625
$ws.Range("J41").Value = 374.5
$ws.Range("K41").Value = 287.625
$ws.Range("L41").Value = 374.5
$ws.Range("M41").Value = 152.375
$ws.Range("N41").Value = -1254.5
$ws.Range("H42").Value = 126.75
$ws.Range("I42").Value = 53.5
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 160.5
$ws.Range("L42").Value = 600
$ws.Range("M42").Value = 69.5
$ws.Range("N42").Value = -1060
$ws.Range("H43").Value = 1586.9722
$ws.Range("I43").Value = 1417.5
$ws.Range("K43").Value = 1417.5
$ws.Range("M43").Value = -1348.5
$ws.Range("H51").Value = 12819.5
$ws.Range("I51").Value = 5373.75
$ws.Range("K51").Value = 5373.75
$ws.Range("M51").Value = -4889.75
$ws.Range("H80").Value = 482.47058
$ws.Range("I80").Value = 636.2857
$ws.Range("J80").Value = 374.8
$ws.Range("K80").Value = 1908.8571
$ws.Range("L80").Value = 1124.4
$ws.Range("M80").Value = -910.8571000000002
$ws.Range("N80").Value = -3120.4
$ws.Range("H82").Value = 6534.75
$ws.Range("I82").Value = 486.2
$ws.Range("K82").Value = 1458.6
$ws.Range("M82").Value = -1052.6
$ws.Range("H83").Value = 482.47058
$ws.Range("I83").Value = 636.2857
$ws.Range("J83").Value = 374.8
$ws.Range("K83").Value = 5726.571300000001
$ws.Range("L83").Value = 3373.2
$ws.Range("M83").Value = -734.5713000000005
$ws.Range("N83").Value = -13357.2
$ws.Range("H85").Value = 6534.75
$ws.Range("I85").Value = 486.2
$ws.Range("K85").Value = 1458.6
$ws.Range("M85").Value = -54.59999999999991
$ws.Range("H88").Value = 4733
$ws.Range("J88").Value = 4733
$ws.Range("L88").Value = 4733
$ws.Range("N88").Value = -5545
$ws.Range("H91").Value = 4733
$ws.Range("J91").Value = 4733
$ws.Range("L91").Value = 4733
$ws.Range("N91").Value = -7541
$ws.Range("H101").Value = 2153.875
$ws.Range("I101").Value = 2850.4
$ws.Range("K101").Value = 8551.200000000001
$ws.Range("M101").Value = -6929.200000000001
$ws.Range("H111").Value = 3620.1
$ws.Range("J111").Value = 10475
$ws.Range("L111").Value = 31425
$ws.Range("N111").Value = -37559
$ws.Range("H113").Value = 3129.1765
$ws.Range("J113").Value = 3179.25
$ws.Range("L113").Value = 3179.25
$ws.Range("N113").Value = -9687.25
$ws.Range("H137").Value = 9219.833000000001
$ws.Range("J137").Value = 12199.2
$ws.Range("L137").Value = 36597.60000000001
$ws.Range("N137").Value = -41697.60000000001
$ws.Range("H138").Value = 3517.2593
$ws.Range("J138").Value = 4200
$ws.Range("L138").Value = 12600
$ws.Range("N138").Value = -22880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3971.0312
$ws.Range("I94").Value = 680
$ws.Range("J94").Value = 9456.083000000001
$ws.Range("K94").Value = 680
$ws.Range("L94").Value = 9456.083000000001
$ws.Range("M94").Value = -229
$ws.Range("N94").Value = -10358.083
$ws.Range("H99").Value = 6652.9473
$ws.Range("I99").Value = 6795.9785
$ws.Range("J99").Value = 5980.7
$ws.Range("K99").Value = 6795.9785
$ws.Range("L99").Value = 5980.7
$ws.Range("M99").Value = -5297.9785
$ws.Range("N99").Value = -8976.700000000001
$ws.Range("H134").Value = 858897.3
$ws.Range("J134").Value = 17531.615
$ws.Range("L134").Value = 52594.845
$ws.Range("N134").Value = -57664.845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3821.4773
$ws.Range("J31").Value = 5858.04
$ws.Range("L31").Value = 5858.04
$ws.Range("N31").Value = -6448.04
$ws.Range("H34").Value = 3821.4773
$ws.Range("J34").Value = 5858.04
$ws.Range("L34").Value = 5858.04
$ws.Range("N34").Value = -6262.04
$ws.Range("H58").Value = 58834172
$ws.Range("I58").Value = 83340856
$ws.Range("K58").Value = 83340856
$ws.Range("M58").Value = -83340653
$ws.Range("H132").Value = 25718.705
$ws.Range("I132").Value = 29001.346
$ws.Range("K132").Value = 87004.038
$ws.Range("M132").Value = -84474.038
$ws.Range("H136").Value = 58834172
$ws.Range("I136").Value = 83340856
$ws.Range("K136").Value = 250022568
$ws.Range("M136").Value = -250020018

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4101.643
$ws.Range("J107").Value = 5229.34
$ws.Range("L107").Value = 15688.02
$ws.Range("N107").Value = -19528.02
$ws.Range("H122").Value = 125040.93
$ws.Range("J122").Value = 160112.56
$ws.Range("L122").Value = 1441013.04
$ws.Range("N122").Value = -1445913.04
$ws.Range("H129").Value = 11905888
$ws.Range("J129").Value = 27779076
$ws.Range("L129").Value = 83337228
$ws.Range("N129").Value = -83347228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3667.0527
$ws.Range("I40").Value = 2346
$ws.Range("K40").Value = 2346
$ws.Range("M40").Value = -2210

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14070.714
$ws.Range("J62").Value = 13436.75
$ws.Range("L62").Value = 13436.75
$ws.Range("N62").Value = -14684.75
$ws.Range("H65").Value = 14070.714
$ws.Range("J65").Value = 13436.75
$ws.Range("L65").Value = 67183.75
$ws.Range("N65").Value = -73423.75
$ws.Range("H100").Value = 1295.125
$ws.Range("I100").Value = 1653
$ws.Range("J100").Value = 937.25
$ws.Range("K100").Value = 3306
$ws.Range("L100").Value = 1874.5
$ws.Range("M100").Value = -2765
$ws.Range("N100").Value = -2956.5
$ws.Range("H132").Value = 21280.75
$ws.Range("I132").Value = 23688.5
$ws.Range("J132").Value = 18873
$ws.Range("K132").Value = 71065.5
$ws.Range("L132").Value = 56619
$ws.Range("M132").Value = -68535.5
$ws.Range("N132").Value = -61679
$ws.Range("H133").Value = 1000000
$ws.Range("J133").Value = 1000000
$ws.Range("L133").Value = 1000000
$ws.Range("N133").Value = -1010120
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360
